$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToShift = @(10, 15, 19, 24, 25)

foreach ($r in $rowsToShift) {
    $ws.Cells.Item($r, 3).Copy()
    $ws.Cells.Item($r - 1, 3).PasteSpecial(-4163)
    $ws.Cells.Item($r, 4).Copy()
    $ws.Cells.Item($r - 1, 4).PasteSpecial(-4163)

    $ws.Cells.Item($r, 3).Value() = ""
    $ws.Cells.Item($r, 4).Value() = ""
}

# The trailing row 28 is now an orphaned empty dialogue row; delete it entirely.
$ws.Rows.Item(28).Delete()

# Clear the stale F13 selection left over from editing.
$ws.Range("A1").Select()
